# Update the "想去人数" (want-to-go count) figures for several events.
# These edits apply to both the "展览" sheet and the "全部类型" sheet,
# which both list the same set of exhibitions.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 72
$wsExhibit.Range("F8").Value  = 109
$wsExhibit.Range("F9").Value  = 8721
$wsExhibit.Range("F12").Value = 1143
$wsExhibit.Range("F13").Value = 982
$wsExhibit.Range("F14").Value = 109
$wsExhibit.Range("F18").Value = 256
$wsExhibit.Range("F19").Value = 66
$wsExhibit.Range("F21").Value = 1027

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 72
$wsAll.Range("F10").Value = 109
$wsAll.Range("F11").Value = 8721
$wsAll.Range("F14").Value = 1143
$wsAll.Range("F15").Value = 982
$wsAll.Range("F16").Value = 109
$wsAll.Range("F20").Value = 256
$wsAll.Range("F21").Value = 66
$wsAll.Range("F23").Value = 1027
